$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# First "Lens" textbox (shape index 14, cNvPr id 17 "TextBox 16") -> "Lens 1"
$lens1 = $s.Shapes.Item(14)
$lens1.TextFrame.TextRange.Text = "Lens 1"

# Second "Lens" textbox (shape index 28, cNvPr id 45 "TextBox 44") -> "Lens 2"
# also reposition/resize it
$lens2 = $s.Shapes.Item(28)
$lens2.TextFrame.TextRange.Text = "Lens 2"
$lens2.Left = 347.3838188976378
$lens2.Top = 161.53311023622047
$lens2.Width = 72.87271653543307
$lens2.Height = 29.081299212598424
